$d = $word.ActiveDocument

# The "Programa" section has two body paragraphs (Portuguese, then the
# italic English translation) that each currently hold one long run of
# text. We split each into four runs separated by manual line breaks
# (<w:br/>) at the sentence/topic boundaries, using Find/Replace with the
# "^l" manual-line-break code so Word inserts a real line break between
# the new runs.

$pt = $d.Paragraphs.Item(14).Range
$pt.Find.Execute("picnometria. Análises microestruturais", $false, $false, $false, $false, $false, $true, 1, $false, "picnometria. ^lAnálises microestruturais", 2)
$pt.Find.Execute("(EDX e WDX). Análises térmicas", $false, $false, $false, $false, $false, $true, 1, $false, "(EDX e WDX). ^lAnálises térmicas", 2)
$pt.Find.Execute("(TGA).Reometria de líquidos", $false, $false, $false, $false, $false, $true, 1, $false, "(TGA).^lReometria de líquidos", 2)

$en = $d.Paragraphs.Item(15).Range
$en.Find.Execute("pycnometry.Microstructural analysis", $false, $false, $false, $false, $false, $true, 1, $false, "pycnometry.^lMicrostructural analysis", 2)
$en.Find.Execute("(EDX and WDX).Thermal analysis", $false, $false, $false, $false, $false, $true, 1, $false, "(EDX and WDX).^lThermal analysis", 2)
$en.Find.Execute("(TGA).Rheometry of liquids", $false, $false, $false, $false, $false, $true, 1, $false, "(TGA).^lRheometry of liquids", 2)
